$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuil1")
$ws2 = $wb.Worksheets.Item("Feuil2")

# --- Sheet1 (Feuil1): new block rows 53-62 ---
$ws1.Range("A53").Value = "Objectif du 22/10/2014"
$ws1.Range("C54").Value = "Ludovic"
$ws1.Range("G55").Value = "Jump mouvement"
$ws1.Range("C57").Value = "Romain"
$ws1.Range("G58").Value = "Design Rope (Vector)"
$ws1.Range("G59").Value = "Hit lvl 1 (Vector) "
$ws1.Range("C61").Value = "Yuxing"
$ws1.Range("G61").Value = "Animation des collectibles (UNITY)"

